$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("B6").Value = "quynhnguyen906@gmail.com"
$ws.Range("B4").Value = "quynhnguyen901@gmail.com"
$ws.Range("C4").Value = "QuynhNguyen"

# Rebuild hyperlinks in the order: B2, B3, B5, B6, B4 (B6 and B4 effectively
# re-created with new targets, landing at the end of the collection)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:quynhnguyen135@gmail.com")
$ws.Range("B2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:quynhnguyen136@gmail.com")
$ws.Range("B3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:quynhnguyen138@gmail.com")
$ws.Range("B5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:quynhnguyen906@gmail.com")
$ws.Range("B6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:quynhnguyen901@gmail.com")
$ws.Range("B4").Style = "Hyperlink"

# Restore selection to B4
$ws.Range("B4").Select()
